$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$xmlTemplate = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>{0}</w:t><w:br/><w:t xml:space="preserve">{1}</w:t><w:br/><w:t xml:space="preserve">{2}</w:t><w:br/><w:t>{3}</w:t><w:br/><w:t>{4}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$cellData = @(
    @{ Row = 1; Col = 1; Header = "69 x 10"; Line2 = "  1    0"; Dashes = "  ----"; Line4 = "6|    |"; Line5 = "9|    |" }
    @{ Row = 1; Col = 2; Header = "21 x 53"; Line2 = "  5    3"; Dashes = "  ----"; Line4 = "2|    |"; Line5 = "1|    |" }
    @{ Row = 1; Col = 3; Header = "58 x 57"; Line2 = "  5    7"; Dashes = "  ----"; Line4 = "5|    |"; Line5 = "8|    |" }
    @{ Row = 2; Col = 1; Header = "14 x 62"; Line2 = "  6    2"; Dashes = "  ----"; Line4 = "1|    |"; Line5 = "4|    |" }
    @{ Row = 2; Col = 2; Header = "22 x 30"; Line2 = "  3    0"; Dashes = "  ----"; Line4 = "2|    |"; Line5 = "2|    |" }
    @{ Row = 2; Col = 3; Header = "60 x 38"; Line2 = "  3    8"; Dashes = "  ----"; Line4 = "6|    |"; Line5 = "0|    |" }
    @{ Row = 3; Col = 1; Header = "98 x 82"; Line2 = "  8    2"; Dashes = "  ----"; Line4 = "9|    |"; Line5 = "8|    |" }
    @{ Row = 3; Col = 2; Header = "22 x 33"; Line2 = "  3    3"; Dashes = "  ----"; Line4 = "2|    |"; Line5 = "2|    |" }
    @{ Row = 3; Col = 3; Header = "89 x 66"; Line2 = "  6    6"; Dashes = "  ----"; Line4 = "8|    |"; Line5 = "9|    |" }
    @{ Row = 4; Col = 1; Header = "28 x 11"; Line2 = "  1    1"; Dashes = "  ----"; Line4 = "2|    |"; Line5 = "8|    |" }
    @{ Row = 4; Col = 2; Header = "19 x 19"; Line2 = "  1    9"; Dashes = "  ----"; Line4 = "1|    |"; Line5 = "9|    |" }
    @{ Row = 4; Col = 3; Header = "54 x 21"; Line2 = "  2    1"; Dashes = "  ----"; Line4 = "5|    |"; Line5 = "4|    |" }
    @{ Row = 5; Col = 1; Header = "96 x 29"; Line2 = "  2    9"; Dashes = "  ----"; Line4 = "9|    |"; Line5 = "6|    |" }
    @{ Row = 5; Col = 2; Header = "79 x 82"; Line2 = "  8    2"; Dashes = "  ----"; Line4 = "7|    |"; Line5 = "9|    |" }
    @{ Row = 5; Col = 3; Header = "62 x 44"; Line2 = "  4    4"; Dashes = "  ----"; Line4 = "6|    |"; Line5 = "2|    |" }
)

foreach ($c in $cellData) {
    $cell = $tbl.Cell($c.Row, $c.Col)
    $xml = $xmlTemplate -f $c.Header, $c.Line2, $c.Dashes, $c.Line4, $c.Line5
    $cell.Range.InsertXML($xml)
}

Write-Host "All 15 cells updated"
